$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 35.995988
$ws.Range("H2").Value = 107.987964
$ws.Range("I2").Value = 0.5613901502831141
$ws.Range("J2").Value = 0.561390150283114
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.044118333333333
$ws.Range("N2").Value = 6.132354999999999
$ws.Range("O2").Value = 0.1776005292722278
$ws.Range("P2").Value = 0.1776005292722278
$ws.Range("Q2").Value = 73.58005899724667
$ws.Range("R2").Value = 662.2205309752198
$ws.Range("S2").Value = 0.09970318781849656
$ws.Range("T2").Value = 0.09970318781849653
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 35.995988
$ws.Range("H3").Value = 107.987964
$ws.Range("I3").Value = 0.5613901502831141
$ws.Range("J3").Value = 0.561390150283114
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.059280333333334
$ws.Range("N3").Value = 21.177841
$ws.Range("O3").Value = 0.6133362746356149
$ws.Range("P3").Value = 0.6133362746356149
$ws.Range("Q3").Value = 254.1057701673027
$ws.Range("R3").Value = 2286.951931505724
$ws.Range("S3").Value = 0.3443209433917732
$ws.Range("T3").Value = 0.3443209433917732
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 35.995988
$ws.Range("H4").Value = 107.987964
$ws.Range("I4").Value = 0.5613901502831141
$ws.Range("J4").Value = 0.561390150283114
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.406242333333333
$ws.Range("N4").Value = 7.218726999999999
$ws.Range("O4").Value = 0.2090631960921573
$ws.Range("P4").Value = 0.2090631960921573
$ws.Range("Q4").Value = 86.61507015575867
$ws.Range("R4").Value = 779.535631401828
$ws.Range("S4").Value = 0.1173660190728443
$ws.Range("T4").Value = 0.1173660190728443
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.666474
$ws.Range("H5").Value = 61.999422
$ws.Range("I5").Value = 0.3223124461726698
$ws.Range("J5").Value = 0.3223124461726698
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.044118333333333
$ws.Range("N5").Value = 6.132354999999999
$ws.Range("O5").Value = 0.1776005292722278
$ws.Range("P5").Value = 0.1776005292722278
$ws.Range("Q5").Value = 42.24471838875665
$ws.Range("R5").Value = 380.2024654988099
$ws.Range("S5").Value = 0.05724286103129258
$ws.Range("T5").Value = 0.05724286103129258
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.666474
$ws.Range("H6").Value = 61.999422
$ws.Range("I6").Value = 0.3223124461726698
$ws.Range("J6").Value = 0.3223124461726698
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.059280333333334
$ws.Range("N6").Value = 21.177841
$ws.Range("O6").Value = 0.6133362746356149
$ws.Range("P6").Value = 0.6133362746356149
$ws.Range("Q6").Value = 145.8904334675447
$ws.Range("R6").Value = 1313.013901207902
$ws.Range("S6").Value = 0.1976859150042375
$ws.Range("T6").Value = 0.1976859150042375
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.666474
$ws.Range("H7").Value = 61.999422
$ws.Range("I7").Value = 0.3223124461726698
$ws.Range("J7").Value = 0.3223124461726698
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.406242333333333
$ws.Range("N7").Value = 7.218726999999999
$ws.Range("O7").Value = 0.2090631960921573
$ws.Range("P7").Value = 0.2090631960921573
$ws.Range("Q7").Value = 49.72854461953266
$ws.Range("R7").Value = 447.5569015757939
$ws.Range("S7").Value = 0.06738367013713975
$ws.Range("T7").Value = 0.06738367013713975
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.456917333333333
$ws.Range("H8").Value = 22.370752
$ws.Range("I8").Value = 0.116297403544216
$ws.Range("J8").Value = 0.116297403544216
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.044118333333333
$ws.Range("N8").Value = 6.132354999999999
$ws.Range("O8").Value = 0.1776005292722278
$ws.Range("P8").Value = 0.1776005292722278
$ws.Range("Q8").Value = 15.24282143121778
$ws.Range("R8").Value = 137.18539288096
$ws.Range("S8").Value = 0.02065448042243863
$ws.Range("T8").Value = 0.02065448042243863
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.456917333333333
$ws.Range("H9").Value = 22.370752
$ws.Range("I9").Value = 0.116297403544216
$ws.Range("J9").Value = 0.116297403544216
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.059280333333334
$ws.Range("N9").Value = 21.177841
$ws.Range("O9").Value = 0.6133362746356149
$ws.Range("P9").Value = 0.6133362746356149
$ws.Range("Q9").Value = 52.64046987849245
$ws.Range("R9").Value = 473.764228906432
$ws.Range("S9").Value = 0.07132941623960422
$ws.Range("T9").Value = 0.07132941623960422
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.456917333333333
$ws.Range("H10").Value = 22.370752
$ws.Range("I10").Value = 0.116297403544216
$ws.Range("J10").Value = 0.116297403544216
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.406242333333333
$ws.Range("N10").Value = 7.218726999999999
$ws.Range("O10").Value = 0.2090631960921573
$ws.Range("P10").Value = 0.2090631960921573
$ws.Range("Q10").Value = 17.94315016363378
$ws.Range("R10").Value = 161.488351472704
$ws.Range("S10").Value = 0.02431350688217318
$ws.Range("T10").Value = 0.02431350688217318
